# Auto-generated edit script: update cryptos list values (Price/Volume columns)
# as produced by the scheduled GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.751.37'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.35%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.805.44'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.93%  '
$ws.Range("E4").Value = '  +0.34%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '232.32'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -1.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5917'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -1.63%  '
$ws.Range("E7").Value = '  +0.35%  '
$ws.Range("E8").Value = '  -0.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06832'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.18%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.35'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.47%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07501'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.90%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.803.94'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -0.79%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.769'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.29%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6241'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.59%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.051.46'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.86%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000009279'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -6.30%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '75.77'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -3.26%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '28.706.25'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.49%  '
$ws.Range("E19").Value = '  -5.95%  '
$ws.Range("E20").Value = '  +0.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '211.39'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -6.32%  '
$ws.Range("E22").Value = '  -1.93%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.842'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -1.66%  '
$ws.Range("E24").Value = '  +0.42%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.30'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.37%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.886'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1270'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.14%  '
$ws.Range("E28").Value = '  -0.44%  '
$ws.Range("E29").Value = '  -4.33%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.06177'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.10%  '
$ws.Range("E31").Value = '  -1.39%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.785'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -0.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.763'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -0.54%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.732'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.065'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -4.86%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.6424'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.08%  '
$ws.Range("E38").Value = '  +0.23%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.592'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.83%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01713'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.29%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.142.07'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -5.73%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8829'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.01%  '
$ws.Range("E43").Value = '  +0.71%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.31'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.962.28'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.83%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '60.53'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000111'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.50%  '
$ws.Range("E48").Value = '  +0.76%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.369'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.25%  '
$ws.Range("E50").Value = '  -0.62%  '
$ws.Range("E51").Value = '  -1.48%  '
